$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.079.92'
$ws.Range("E2").Value = '  -3.38%  '
$ws.Range("D3").Value = '3.721.30'
$ws.Range("E3").Value = '  -4.10%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '616.16'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.88'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +4.82%  '
$ws.Range("D7").Value = '3.715.62'
$ws.Range("E7").Value = '  -4.13%  '
$ws.Range("E8").Value = '  -5.33%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("E10").Value = '  -3.67%  '
$ws.Range("E11").Value = '  -8.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '57.59'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +6.18%  '
$ws.Range("E13").Value = '  -8.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.74'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -6.11%  '
$ws.Range("D15").Value = '4.308.68'
$ws.Range("E15").Value = '  -4.14%  '
$ws.Range("D16").Value = '3.719.71'
$ws.Range("E16").Value = '  -4.31%  '
$ws.Range("E17").Value = '  -7.15%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.01'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -6.92%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.126'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.14'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -6.82%  '
$ws.Range("D21").Value = '68.791.99'
$ws.Range("E21").Value = '  -3.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '415.68'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -5.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.74'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.99%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '89.32'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -5.21%  '
$ws.Range("E25").Value = '  -7.85%  '
$ws.Range("E26").Value = '  -7.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.00'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -6.55%  '
$ws.Range("E28").Value = '  -3.62%  '
$ws.Range("E29").Value = '  +1.74%  '
$ws.Range("E30").Value = '  -8.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.13'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -5.92%  '
$ws.Range("E32").Value = '  -15.92%  '
$ws.Range("E33").Value = '  -7.61%  '
$ws.Range("E34").Value = '  -5.74%  '
$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '65.71'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.87%  '
$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '44.07'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -8.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '610.56'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.25%  '
$ws.Range("E38").Value = '  -11.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.408'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -6.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.21%  '
$ws.Range("E42").Value = '  -5.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.06'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -7.93%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0443'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -6.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.67'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -7.62%  '
$ws.Range("B46").Value = 'THORChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.27'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -9.46%  '
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.79'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -12.82%  '
$ws.Range("E48").Value = '  -6.11%  '
$ws.Range("D49").Value = '2.797.56'
$ws.Range("E49").Value = '  -3.97%  '
$ws.Range("E50").Value = '  -7.04%  '
$ws.Range("B51").Value = 'FLOKI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000265'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -5.52%  '
